$d = $word.ActiveDocument

# --- Change 1: merge the split runs in the "Social Media Links" paragraph
# into a single run (simple clean-up / no wording change).
$d.Content.Find.Execute(
    "by " + "clicking" + " on a logo below to take you to the appropriate " + "website" + " or by pressing",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by clicking on a logo below to take you to the appropriate website or by pressing",
    2)

# --- Change 2: insert a new blank paragraph + a new "This is Gill!" canvas
# paragraph right after the CV paragraph (and before the existing blank
# paragraph that follows it).
$rng = $d.Content
$rng.Find.Execute("continue to look around")
$rng.Collapse(0)
$rng.InsertParagraphAfter()

$newPara = $d.Paragraphs(19)

$gillPara = '<w:p><w:r><w:t xml:space="preserve">This is Gill! As you can see, he is untextured </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>at the moment</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">, but you can still view his model in all his glory while he stands in a T-Pose for you. Please feel free to interact with the viewer using your mouse, </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>clicking</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and dragging to spin him, and using your scroll wheel to zoom in and out, or pinch and pull if you are using a trackpad/touchscreen device. </w:t></w:r></w:p>'

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/>' + $gillPara + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xml)
